$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: rows 4-7 (499659a2, 7674aa93, 8cd62b24, b0a124d1) ---
# Priority (col E) moves from "low" to "ht", and Latest Handoff Datetime
# (col H) is refreshed to the newly-generated handoff timestamp.
$wsZh = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $wsZh.Cells.Item($r, 5).Value = "ht"
    $wsZh.Cells.Item($r, 8).Value = "2016-09-08 04:44:04"
}

# --- de-de sheet: rows 4-7 (499659a2, 7674aa93, 8cd62b24, b0a124d1) ---
# Same Priority change; its Latest Handoff Datetime column shares the
# "Latest HO Xliff Generate Date" string with the Overview sheet.
$wsDe = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $wsDe.Cells.Item($r, 5).Value = "ht"
    $wsDe.Cells.Item($r, 8).Value = "2016-09-08 04:44:14"
}

# --- Overview sheet: rows 4-7, column G "Latest HO Xliff Generate Date" ---
# Shares the same shared-string slot as de-de!H4:H7, so it updates in lockstep.
$wsOverview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-09-08 04:44:14"
}
